$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = 16.613
$ws.Range("B8").Value = 6.185
$ws.Range("B10").Value = 6.188000000000001
$ws.Range("B12").Value = 5.315
$ws.Range("E15").Value = 16.361
$ws.Range("B18").Value = 5.137
$ws.Range("E18").Value = 16.466
$ws.Range("E20").Value = 16.406
$ws.Range("E29").Value = 17.09
$ws.Range("E30").Value = 16.37
$ws.Range("E31").Value = 16.426
$ws.Range("B37").Value = 8.73
$ws.Range("E40").Value = 16.627
$ws.Range("E50").Value = 16.291
$ws.Range("B55").Value = 4.572
$ws.Range("B68").Value = 5.220000000000001
$ws.Range("E68").Value = 17.272
$ws.Range("E76").Value = 16.623
$ws.Range("B77").Value = 6.123
$ws.Range("B78").Value = 7.810999999999998
$ws.Range("B81").Value = 6.008
$ws.Range("B82").Value = 5.658999999999999
$ws.Range("E87").Value = 16.272
$ws.Range("E88").Value = 16.224
$ws.Range("E96").Value = 16.38
$ws.Range("E98").Value = 16.282
$ws.Range("E101").Value = 16.716
$ws.Range("E102").Value = 16.649
